$d = $word.ActiveDocument

# The trailing "_GoBack" bookmark currently sits at the end of the 페이징
# paragraph. It needs to move to the end of the document (end of the new
# "스레싱" paragraph that this edit appends), so delete it here and it is
# re-created in the right spot as part of the appended XML below.
$d.Bookmarks.ShowHidden = $true
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$d.Bookmarks.ShowHidden = $false

# Paragraph 33 ("세스먼테이션 기법 ...") gets its paragraph-mark formatting
# simplified (no explicit eastAsia rFonts on the mark) and its two adjoining
# descriptive runs merged into a single run.
$segPara = $d.Paragraphs.Item(33)
$segPara.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>세스먼테이션</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> 기법 </w:t></w:r><w:r><w:t xml:space="preserve">- </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">가상기억장치에 보관되어 있는 프로그램을 다양한 크기의 논리적인 단위로 나눈 후 주기억장치에 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>적재시켜</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> 실행시키는 기법 </w:t></w:r></w:p>')

# Append the new "136 페이지 교체 알고리즘" / "137 가상기억장치 기타 관리 사항"
# material as a single block of raw paragraph XML at the very end of the
# document body (immediately before sectPr). Using one InsertXML call (rather
# than one per paragraph) keeps the genuinely-empty spacer paragraphs intact.
$endPos = $d.Content.End
$endRange = $d.Range($endPos, $endPos)
$endRange.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p><w:pPr><w:rPr><w:sz w:val="22"/><w:shd w:val="pct15" w:color="auto" w:fill="FFFFFF"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="22"/><w:shd w:val="pct15" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">136 </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:shd w:val="pct15" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>페이지 교체 알고리즘</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>- 페이지 부재가 발생하면 어떤 페이지 프레임을 선택하여 교체할 것인지를 결정하는 기법</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>OTP(최적 교체)</w:t></w:r><w:r><w:t xml:space="preserve"> – </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>앞으로 가장 오랫동안 사용하지 않을 페이지를 교체하는 방법</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>FIFO</w:t></w:r><w:r><w:t xml:space="preserve"> – </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>가장 먼저 들어와서 가장 오래 있었던 페이지를 교체하는 기법</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">LRU </w:t></w:r><w:r><w:t>–</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> 최근에 가장 오랫동안 사용하지 않은 페이지를 교체하는 기법</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">LFU </w:t></w:r><w:r><w:t>–</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> 사용 빈도가 가장 적은 페이지를 교체하는 기법</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">NUR </w:t></w:r><w:r><w:t>–</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> 최근에 사용하지 않은 페이지를 교체하는 기법</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">SCR </w:t></w:r><w:r><w:t>–</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> 가장 오랫동안 주기억장치에 있던 페이지 중 자주 사용되는 페이지의 교체를 방지하기 위한 기법</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:rPr><w:sz w:val="22"/><w:shd w:val="pct15" w:color="auto" w:fill="FFFFFF"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:sz w:val="22"/><w:shd w:val="pct15" w:color="auto" w:fill="FFFFFF"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:sz w:val="22"/><w:shd w:val="pct15" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:lastRenderedPageBreak/><w:t>137 가상기억장치 기타 관리 사항</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">Locality </w:t></w:r><w:r><w:t>–</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> 프로세스가 실행되는 동안 주기억장치를 참조할 때 일부 페이지만 집중적으로 참조하는 성질이 있음</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">시간 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>구역성</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">/ </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">공간 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>구역성</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">워킹 셋 </w:t></w:r><w:r><w:t xml:space="preserve">– </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>프로세스가 일정시간 동안 자주 참조하는 페이지들의 집합</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>스레싱</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">– </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve">프로세스의 처리 시간보다 페이지 교체에 소용되는 시간이 많아지는 현상 </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>')
